$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 105.9375
$ws.Range("I9").Value = 56.333332
$ws.Range("J9").Value = 254.75
$ws.Range("K9").Value = 56.333332
$ws.Range("L9").Value = 254.75
$ws.Range("M9").Value = 112.666668
$ws.Range("N9").Value = -592.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7814031.5
$ws.Range("I19").Value = 15625225
$ws.Range("J19").Value = 2837.5
$ws.Range("K19").Value = 15625225
$ws.Range("L19").Value = 2837.5
$ws.Range("M19").Value = -15625050
$ws.Range("N19").Value = -3187.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 57611.4
$ws.Range("I21").Value = 57611.4
$ws.Range("K21").Value = 57611.4
$ws.Range("M21").Value = -57143.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 57611.4
$ws.Range("I23").Value = 57611.4
$ws.Range("K23").Value = 57611.4
$ws.Range("M23").Value = -57377.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 839.1818
$ws.Range("I32").Value = 753.44446
$ws.Range("J32").Value = 1225
$ws.Range("K32").Value = 753.44446
$ws.Range("L32").Value = 1225
$ws.Range("M32").Value = -427.44446
$ws.Range("N32").Value = -1877

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1412.7778
$ws.Range("I40").Value = 1816.6666
$ws.Range("J40").Value = 1210.8334
$ws.Range("K40").Value = 1816.6666
$ws.Range("L40").Value = 1210.8334
$ws.Range("M40").Value = -1641.6666
$ws.Range("N40").Value = -1560.8334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2413
$ws.Range("I43").Value = 1001
$ws.Range("J43").Value = 2695.4
$ws.Range("K43").Value = 1001
$ws.Range("L43").Value = 2695.4
$ws.Range("M43").Value = -932
$ws.Range("N43").Value = -2833.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 11776.538
$ws.Range("I53").Value = 30166
$ws.Range("J53").Value = 283.125
$ws.Range("K53").Value = 30166
$ws.Range("L53").Value = 283.125
$ws.Range("M53").Value = -29529
$ws.Range("N53").Value = -1557.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 32341586
$ws.Range("I98").Value = 14286210
$ws.Range("J98").Value = 63938492
$ws.Range("K98").Value = 14286210
$ws.Range("L98").Value = 63938492
$ws.Range("M98").Value = -14284712
$ws.Range("N98").Value = -63941488

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4350778.5
$ws.Range("I113").Value = 12502112
$ws.Range("J113").Value = 3400
$ws.Range("K113").Value = 12502112
$ws.Range("L113").Value = 3400
$ws.Range("M113").Value = -12498858
$ws.Range("N113").Value = -9908

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 16675089
$ws.Range("I116").Value = 8335696
$ws.Range("J116").Value = 27794278
$ws.Range("K116").Value = 8335696
$ws.Range("L116").Value = 27794278
$ws.Range("M116").Value = -8332254
$ws.Range("N116").Value = -27801162

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 32341586
$ws.Range("I122").Value = 14286210
$ws.Range("J122").Value = 63938492
$ws.Range("K122").Value = 42858630
$ws.Range("L122").Value = 191815476
$ws.Range("M122").Value = -42856180
$ws.Range("N122").Value = -191820376

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2780425.8
$ws.Range("I132").Value = 697541.4399999999
$ws.Range("K132").Value = 2092624.32
$ws.Range("M132").Value = -2090094.32

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2326.75
$ws.Range("I141").Value = 1749.3334
$ws.Range("K141").Value = 5248.0002
$ws.Range("M141").Value = -68.0002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2540024.2
$ws.Range("I32").Value = 3181591
$ws.Range("J32").Value = 13855
$ws.Range("K32").Value = 3181591
$ws.Range("L32").Value = 13855
$ws.Range("M32").Value = -3181304
$ws.Range("N32").Value = -14429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3046188
$ws.Range("I61").Value = 1345134.8
$ws.Range("J61").Value = 29412514
$ws.Range("K61").Value = 1345134.8
$ws.Range("L61").Value = 29412514
$ws.Range("M61").Value = -1344922.8
$ws.Range("N61").Value = -29412938

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 45098920
$ws.Range("I74").Value = 34483308
$ws.Range("J74").Value = 106669470
$ws.Range("K74").Value = 34483308
$ws.Range("L74").Value = 106669470
$ws.Range("M74").Value = -34482434
$ws.Range("N74").Value = -106671218

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 45098920
$ws.Range("I77").Value = 34483308
$ws.Range("J77").Value = 106669470
$ws.Range("K77").Value = 172416540
$ws.Range("L77").Value = 533347350
$ws.Range("M77").Value = -172412172
$ws.Range("N77").Value = -533356086

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2471.5293
$ws.Range("I110").Value = 1400.3
$ws.Range("J110").Value = 4001.8572
$ws.Range("K110").Value = 1400.3
$ws.Range("L110").Value = 4001.8572
$ws.Range("M110").Value = 644.7
$ws.Range("N110").Value = -8091.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 49788.535
$ws.Range("J135").Value = 49788.535
$ws.Range("L135").Value = 49788.535
$ws.Range("N135").Value = -59928.535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3046188
$ws.Range("I136").Value = 1345134.8
$ws.Range("J136").Value = 29412514
$ws.Range("K136").Value = 4035404.4
$ws.Range("L136").Value = 88237542
$ws.Range("M136").Value = -4032854.4
$ws.Range("N136").Value = -88242642

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 33337696
$ws.Range("I64").Value = 506
$ws.Range("J64").Value = 41671990
$ws.Range("K64").Value = 506
$ws.Range("L64").Value = 41671990
$ws.Range("M64").Value = -281
$ws.Range("N64").Value = -41672440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 33337696
$ws.Range("I67").Value = 506
$ws.Range("J67").Value = 41671990
$ws.Range("K67").Value = 506
$ws.Range("L67").Value = 41671990
$ws.Range("M67").Value = 274
$ws.Range("N67").Value = -41673550

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 335.46667
$ws.Range("I80").Value = 391.5
$ws.Range("J80").Value = 315.0909
$ws.Range("K80").Value = 391.5
$ws.Range("L80").Value = 315.0909
$ws.Range("M80").Value = 606.5
$ws.Range("N80").Value = -2311.0909

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 335.46667
$ws.Range("I83").Value = 391.5
$ws.Range("J83").Value = 315.0909
$ws.Range("K83").Value = 1957.5
$ws.Range("L83").Value = 1575.4545
$ws.Range("M83").Value = 3034.5
$ws.Range("N83").Value = -11559.4545

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1115.8334
$ws.Range("I94").Value = 806.9048
$ws.Range("J94").Value = 1836.6666
$ws.Range("K94").Value = 806.9048
$ws.Range("L94").Value = 1836.6666
$ws.Range("M94").Value = -355.9048
$ws.Range("N94").Value = -2738.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 13937487
$ws.Range("I134").Value = 15205458
$ws.Range("J134").Value = 6752316.5
$ws.Range("K134").Value = 45616374
$ws.Range("L134").Value = 20256949.5
$ws.Range("M134").Value = -45613839
$ws.Range("N134").Value = -20262019.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 45796
$ws.Range("J135").Value = 45796
$ws.Range("L135").Value = 45796
$ws.Range("N135").Value = -55936

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 835
$ws.Range("I22").Value = 258.66666
$ws.Range("J22").Value = 3716.6667
$ws.Range("K22").Value = 258.66666
$ws.Range("L22").Value = 3716.6667
$ws.Range("M22").Value = 91.33334000000002
$ws.Range("N22").Value = -4416.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 33336848
$ws.Range("I94").Value = 1389.3334
$ws.Range("J94").Value = 55560490
$ws.Range("K94").Value = 1389.3334
$ws.Range("L94").Value = 55560490
$ws.Range("M94").Value = -938.3334
$ws.Range("N94").Value = -55561392

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3363.7222
$ws.Range("I122").Value = 4734.5
$ws.Range("J122").Value = 622.1667
$ws.Range("K122").Value = 14203.5
$ws.Range("L122").Value = 1866.5001
$ws.Range("M122").Value = -11753.5
$ws.Range("N122").Value = -6766.5001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2382.7827
$ws.Range("I132").Value = 1622.4615
$ws.Range("K132").Value = 4867.3845
$ws.Range("M132").Value = -2337.3845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1057581.4
$ws.Range("I134").Value = 4978.5864
$ws.Range("K134").Value = 14935.7592
$ws.Range("M134").Value = -12400.7592

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2484.7273
$ws.Range("J132").Value = 3420.2856
$ws.Range("L132").Value = 30782.5704
$ws.Range("N132").Value = -35842.5704

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 1000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 1000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -1344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 879.825
$ws.Range("I46").Value = 808.3714
$ws.Range("J46").Value = 1380
$ws.Range("K46").Value = 808.3714
$ws.Range("L46").Value = 1380
$ws.Range("M46").Value = -620.3714
$ws.Range("N46").Value = -1756
